# The deck ships two theme parts:
#   theme/theme1.xml -> "Office Theme" (default blue "Office" palette)
#   theme/theme2.xml -> "Integral" (the "Red Violet" palette), which is the
#                        theme actually wired to the one slide master used
#                        by every slide in the deck.
#
# The authored edit swaps the content of those two theme parts: the master's
# theme (theme2.xml) becomes the plain "Office Theme"/"Office" palette, and
# the previously-unused theme1.xml becomes "Integral"/"Red Violet".
#
# Re-theme the presentation's master (and therefore every slide, since they
# all share the one master/design) from "Integral" back to the stock
# "Office Theme" palette, using the Design/ColorScheme object model.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

$colorScheme = $master.ColorScheme

# ppColorScheme index -> target Office-theme RGB (VBA RGB() packs as
# R + G*256 + B*65536, i.e. the same order python-pptx / OOXML srgbClr uses
# but little-endian across the three bytes).
$officeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388     # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501      # accent2  ED7D31
    7  = 10855845     # accent3  A5A5A5
    8  = 49407        # accent4  FFC000
    9  = 12874308     # accent5  4472C4
    10 = 4697456      # accent6  70AD47
    11 = 12673797     # hlink    0563C1
    12 = 7491477       # folHlink 954F72
}

foreach ($index in 1..12) {
    $colorScheme.Colors($index).RGB = $officeColors[$index]
}
